$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.772.03"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "2.583.43"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").Value = "3.044.94"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "62.712.08"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").Value = "2.588.77"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "338.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.69%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("E25").Value = "  +3.75%  "
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("E27").Value = "  -2.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").Value = "0.0₃0811"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "454.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.74%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "176.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.401"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "160.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.01%  "
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.631"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0536"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0964"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("E51").Value = "  -1.71%  "
